# Auto-generated edit script: updates crypto price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '26.988.90' }
    @{ Cell = "E2"; Value = '  +0.26%  ' }
    @{ Cell = "D3"; Value = '1.671.29' }
    @{ Cell = "E3"; Value = '  +0.17%  ' }
    @{ Cell = "D5"; Value = '214.72' }
    @{ Cell = "E5"; Value = '  -0.25%  ' }
    @{ Cell = "D6"; Value = '0.516' }
    @{ Cell = "E6"; Value = '  -0.79%  ' }
    @{ Cell = "E7"; Value = '  +0.04%  ' }
    @{ Cell = "E8"; Value = '  +1.67%  ' }
    @{ Cell = "D9"; Value = '21.40' }
    @{ Cell = "E9"; Value = '  +5.29%  ' }
    @{ Cell = "E10"; Value = '  -0.11%  ' }
    @{ Cell = "D11"; Value = '0.0887' }
    @{ Cell = "E11"; Value = '  -0.44%  ' }
    @{ Cell = "D12"; Value = '1.908.35' }
    @{ Cell = "E12"; Value = '  +0.25%  ' }
    @{ Cell = "D13"; Value = '1.693.88' }
    @{ Cell = "E13"; Value = '  +1.43%  ' }
    @{ Cell = "E14"; Value = '  +0.74%  ' }
    @{ Cell = "D15"; Value = '0.533' }
    @{ Cell = "E15"; Value = '  +1.52%  ' }
    @{ Cell = "D16"; Value = '66.13' }
    @{ Cell = "E16"; Value = '  +0.75%  ' }
    @{ Cell = "B17"; Value = 'WrappedBTC' }
    @{ Cell = "C17"; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' }
    @{ Cell = "D17"; Value = '26.994.19' }
    @{ Cell = "E17"; Value = '  +0.31%  ' }
    @{ Cell = "B18"; Value = 'Chainlink' }
    @{ Cell = "C18"; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link' }
    @{ Cell = "D18"; Value = '8.17' }
    @{ Cell = "E18"; Value = '  +2.66%  ' }
    @{ Cell = "D19"; Value = '234.73' }
    @{ Cell = "E19"; Value = '  -0.28%  ' }
    @{ Cell = "D20"; Value = '0.0₃0734' }
    @{ Cell = "E20"; Value = '  +0.32%  ' }
    @{ Cell = "E21"; Value = '  -0.03%  ' }
    @{ Cell = "D22"; Value = '4.46' }
    @{ Cell = "E22"; Value = '  +1.39%  ' }
    @{ Cell = "D23"; Value = '9.24' }
    @{ Cell = "E23"; Value = '  +0.67%  ' }
    @{ Cell = "E24"; Value = '  -2.59%  ' }
    @{ Cell = "E25"; Value = '  +0.25%  ' }
    @{ Cell = "D26"; Value = '7.24' }
    @{ Cell = "E26"; Value = '  +1.82%  ' }
    @{ Cell = "D27"; Value = '16.37' }
    @{ Cell = "E27"; Value = '  +3.00%  ' }
    @{ Cell = "E28"; Value = '  -0.24%  ' }
    @{ Cell = "E29"; Value = '  -0.05%  ' }
    @{ Cell = "E30"; Value = '  +0.50%  ' }
    @{ Cell = "D32"; Value = '3.36' }
    @{ Cell = "E32"; Value = '  +0.62%  ' }
    @{ Cell = "D33"; Value = '1.536.06' }
    @{ Cell = "E33"; Value = '  +6.12%  ' }
    @{ Cell = "D34"; Value = '3.16' }
    @{ Cell = "E34"; Value = '  +0.58%  ' }
    @{ Cell = "D35"; Value = '1.70' }
    @{ Cell = "E35"; Value = '  +3.89%  ' }
    @{ Cell = "E36"; Value = '  -1.20%  ' }
    @{ Cell = "E37"; Value = '  +0.93%  ' }
    @{ Cell = "E38"; Value = '  +2.17%  ' }
    @{ Cell = "D39"; Value = '0.909' }
    @{ Cell = "E39"; Value = '  +0.61%  ' }
    @{ Cell = "E40"; Value = '  +4.56%  ' }
    @{ Cell = "E41"; Value = '  -0.02%  ' }
    @{ Cell = "D42"; Value = '67.52' }
    @{ Cell = "E42"; Value = '  +2.26%  ' }
    @{ Cell = "D43"; Value = '5.52' }
    @{ Cell = "E43"; Value = '  -3.66%  ' }
    @{ Cell = "E44"; Value = '  -2.56%  ' }
    @{ Cell = "D45"; Value = '1.815.69' }
    @{ Cell = "E45"; Value = '  +0.59%  ' }
    @{ Cell = "D46"; Value = '0.779' }
    @{ Cell = "E46"; Value = '  -0.17%  ' }
    @{ Cell = "D47"; Value = '90.41' }
    @{ Cell = "E47"; Value = '  -0.40%  ' }
    @{ Cell = "E48"; Value = '  +0.40%  ' }
    @{ Cell = "E49"; Value = '  -0.26%  ' }
    @{ Cell = "D50"; Value = '0.103' }
    @{ Cell = "E50"; Value = '  +1.72%  ' }
    @{ Cell = "D51"; Value = '8.02' }
    @{ Cell = "E51"; Value = '  +5.79%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text interpretation so numeric-looking strings (e.g. "21.40")
    # keep their exact formatting instead of being coerced to a Double.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    # Reset to the default (unstyled) cell style so we do not leave a
    # stray explicit number-format behind on cells that originally had none.
    $rng.Style = "Normal"
}
